# DataConfig.xlsx edit — "feat: fix database for juli"
#
# The sheet lists a channel's group username / channel name / post text.
# The old "polkercrew" / "POLKER Crew" row is replaced by a
# "satoshistreetbets" / "SatoshiStreetBets" row, and the Post text (Julio's
# TipTok pitch) gets a small typo fix ("webiste" -> "website") while keeping
# the same post text in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GroupUserName (A2) and ChannelName (B2)
$ws.Range("A2").Value = "satoshistreetbets"
$ws.Range("B2").Value = "SatoshiStreetBets"

# Post (C2) — same bio text, with "webiste" corrected to "website"
$ws.Range("C2").Value = 'Hi, my name is Julio Qose and I''m the CEO of TipTok 😁 TipTok aims to become the Cryptocurrency for Social Media and the next $1 Billion Market Cap token. We have regular AMAs and the whole team is doxxed in the website. For more information on our fascinating roadmap and our upcoming Whitelist Presale visit our website at https://tiptok.finance/ and join our TG https://t.me/TipTok_Official. Join our movement and be part of something great.'

# Selection moved from C10 to C5
$ws.Range("C5").Select()
